# Update the "Synthetic Control" (column B) estimates for each sector sheet
# and strip the bold/bordered header style from column A, to match the
# corrected model output, then add a new "Inference" sheet with RMSPE
# diagnostics for each sector's synthetic control fit.

$wb = $excel.ActiveWorkbook

$sectorData = @{
    "Construction" = @(4.909354875538632, 4.681860772777387, 4.791049753811859, 5.139543440038089, 5.34886721368235, 6.205922490123703, 5.784947119907189, 6.143480909695239, 6.400035322846863, 6.191352450004326, 6.191166133080032, 5.841801154118497, 5.477257431048647, 5.72953994413719, 5.522285612648909, 5.759116819788894, 5.943320960741827, 6.179500522456502, 6.530929923641208, 6.794157231708999, 7.362478398071418, 7.586174354828828, 7.246887571771938)
    "Other_Sectors" = @(70.03696071999781, 71.03745231059301, 71.56292968461145, 70.36612088251015, 70.96063091353383, 69.91490152871764, 69.16778917520644, 67.9151126453994, 67.25923112711487, 66.91384245710941, 67.6048921826494, 67.02533444895845, 67.10089822192366, 66.63574525273724, 67.86074708678451, 67.31400499545701, 67.0483170069154, 67.01120955811074, 66.39476429267414, 66.37016209134727, 66.5310739953849, 67.53666946896473, 67.94759804405861)
    "Financial" = @(3.00992281240767, 3.119219188670747, 3.20044566148106, 3.350867203753796, 3.581513195913283, 3.611754875812044, 3.592473108080242, 3.427637427889325, 3.287682421441319, 3.670919484560697, 4.311852415028941, 3.718027642844674, 3.892211620548045, 3.874455470866506, 3.939770242643672, 4.216545727569942, 3.968085253486091, 3.985776811500933, 4.219435331686264, 4.137767791437408, 3.979826737113193, 4.19878336860629, 4.290798380046681)
    "Agriculture" = @(5.318282310693759, 6.80438705138756, 7.126512227482343, 5.307523136718197, 5.719130016857253, 5.249898607650811, 4.541407462255786, 4.63281278389741, 5.524760988745911, 4.495523460456369, 5.094819016888716, 5.760159237561744, 5.699048411927259, 4.862708248163462, 6.403922498250678, 4.464821979233463, 4.101512155680004, 5.032137675800048, 5.244192295436321, 4.882554067475919, 5.034262621290559, 5.015181403765334, 5.100594466202272)
    "Manufacturing" = @(15.53371666517171, 15.53530257137981, 15.11131679546708, 15.236606742019, 14.04929705430803, 13.71605510353847, 13.3846279909916, 12.8659644013366, 12.45645853700634, 11.90194394046013, 10.89402385347335, 10.96630295927072, 10.95956756856186, 10.55122509245605, 10.08793780686504, 10.73871690942542, 10.92449051656184, 10.05351042074362, 10.03859794262687, 9.696276141096854, 9.395790122008396, 9.025266848704593, 8.752646223152736)
    "Rental," = @(5.621317299995464, 5.564318811431222, 5.467196411875582, 5.799167296374685, 6.017210646691405, 6.198518981985159, 6.378247598254934, 6.523318826965772, 6.176825746001917, 6.689039777759969, 7.261371899193909, 7.04854462044021, 7.197291111761587, 7.25880297214124, 7.109814302637266, 6.974612893596587, 7.11063907522167, 7.284064222783096, 6.964782526682914, 7.173406836040595, 7.368357866644308, 7.64478579410545, 7.372003563562762)
}

foreach ($sheetName in $sectorData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $values = $sectorData[$sheetName]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 2).Value = $values[$i]
    }
}

# Remove the special bold+bordered style that used to mark the "Year" column
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.ClearFormats()
}

# Add the new "Inference" sheet at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$infSheet = $wb.Worksheets.Add($null, $lastSheet)
$infSheet.Name = "Inference"

$headers = @("Sector", "RMSPE_Pre", "RMSPE_Post", "RMSPE_Ratio", "Mean_Post_Gap")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $infSheet.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$inferenceRows = @(
    @("Construction", 0.1568176074244797, 2.499739061241317, 15.94042341479507, 2.180396810056857),
    @("Other_Sectors", 5.531345688477642, 3.66561950609198, 0.6626994067154109, 3.341445970545716),
    @("Financial and Insurance Services", 0.1568777600068354, 0.9196888216176277, 5.862455083356336, -0.8637212587990629),
    @("Agriculture", 0.6170032752778043, 1.524037215625142, 2.470063412449387, 1.375814103515723),
    @("Manufacturing", 1.321982956458023, 1.147922369255937, 0.8683337131149976, 1.103683969312528),
    @("Rental, Hiring and Real Estate Services", 0.1847312011354326, 0.4340397431970666, 2.349574628050285, -0.3173077074809139)
)

for ($r = 0; $r -lt $inferenceRows.Length; $r++) {
    $rowData = $inferenceRows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $infSheet.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}

# Restore the originally-active tab (Construction) rather than leaving focus
# on the newly added sheet
$wb.Worksheets.Item("Construction").Activate()

Write-Host "Inference sheet added; sector values updated; header styling cleared."
